$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for product "SPCT012" (duplicate data of the SPCT011 row above it) ---
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).RowHeight = 20

$ws.Range("A7").Value = 3.0
$ws.Range("B7").Value = 'SPCT012'
$ws.Range("C7").Value = 'AH-X9XEW1'
$ws.Range("D7").Value = 32.0
$ws.Range("E7").Value = 141200.0
$ws.Range("F7").Value = 54300.0
$ws.Range("G7").Value = 'Phổ biến'
$ws.Range("H7").Value = 'Nâu'
$ws.Range("I7").Value = 20.0
$ws.Range("J7").Value = '82.0 - 20.5 - 48.2'
$ws.Range("K7").Value = 'Sắt không gỉ'
$ws.Range("L7").Value = 'maylanh.png'
$ws.Range("M7").Value = 'No Mô Tả'
$ws.Range("N7").Value = 'Đang kinh doanh'

# --- Renumber the STT ("A") column for every row that shifted down one position ---
$ws.Range("A8").Value = 4.0
$ws.Range("A9").Value = 5.0
$ws.Range("A10").Value = 6.0
$ws.Range("A11").Value = 7.0
$ws.Range("A12").Value = 8.0
$ws.Range("A13").Value = 9.0
$ws.Range("A14").Value = 10.0
$ws.Range("A15").Value = 11.0

# --- Update product SPCT09 (now row 15) with its corrected attributes ---
$ws.Range("H15").Value = 'Cà rốt'
$ws.Range("J15").Value = '11.0 - 33.0 - 22.0'
$ws.Range("K15").Value = 'Đá sức mạnh'

# --- Append two brand-new products: SPCT11 and SPCT13 ---
$ws.Rows.Item(16).RowHeight = 20
$ws.Rows.Item(17).RowHeight = 20

$ws.Range("A16").Value = 12.0
$ws.Range("B16").Value = 'SPCT11'
$ws.Range("C16").Value = 'Inverter 11'
$ws.Range("D16").Value = 11.0
$ws.Range("E16").Value = 1100.0
$ws.Range("F16").Value = 31100.0
$ws.Range("G16").Value = 'Không phổ biến'
$ws.Range("H16").Value = 'Nâu'
$ws.Range("I16").Value = 7.0
$ws.Range("J16").Value = '11.0 - 33.0 - 22.0'
$ws.Range("K16").Value = 'Đá sức mạnh'
$ws.Range("L16").Value = 'loa.png'
$ws.Range("M16").Value = 'No Mô Tả'
$ws.Range("N16").Value = 'Đang kinh doanh'

$ws.Range("A17").Value = 13.0
$ws.Range("B17").Value = 'SPCT13'
$ws.Range("C17").Value = 'BM01'
$ws.Range("D17").Value = 1000.0
$ws.Range("E17").Value = 2500.0
$ws.Range("F17").Value = 4000.0
$ws.Range("G17").Value = 'Không phổ biến'
$ws.Range("H17").Value = 'Đen'
$ws.Range("I17").Value = 2.0
$ws.Range("J17").Value = '30.0 - 15.5 - 14.0'
$ws.Range("K17").Value = 'Đá thời gian'
$ws.Range("L17").Value = 'loa.png'
$ws.Range("M17").Value = 'No Mô Tả'
$ws.Range("N17").Value = 'Đang kinh doanh'
